# [10-05] Danh gia nhom 17.
# Update the evaluation header labels for "lan 3/4/5" to include their
# dates, and record the group-17 (rows 11-15) "lan 2 (10/05)" scores.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row 5: give "DG lan 3/4/5" their date suffixes ------------
# Write the new text first, then re-apply the original header formatting
# (wrap text + quote-prefixed left/top alignment) by copying the format
# from a sibling header cell that already carries it, so the style index
# used by H5/I5/J5 matches the "wrapped" header style instead of drifting
# to a freshly-synthesized one.
$ws.Range("H5").Value = "ĐG lần 3" + [char]10 + "(13/05)"
$ws.Range("G5").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null

$ws.Range("I5").Value = "ĐG lần 4" + [char]10 + "(17/05)"
$ws.Range("G5").Copy() | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null

$ws.Range("J5").Value = "ĐG lần 5" + [char]10 + "(20/05)"
$ws.Range("G5").Copy() | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Group 17 (rows 11-15) scored for "DG lan 2 (10/05)" -> column G --
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1

# --- Leave the view the way the author left it: scrolled back to the
# top, with the cursor resting on J6. -----------------------------------
$ws.Range("J6").Select() | Out-Null
